# Enhanced EDP bulk upload functionality with new changes.
# Update the sample data with a new header row (and related view tweaks).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: swap in the new header text / column order ---------------
# Old header: EDP | Material Description | Section | Category | UoM
# New header: Material | Material Description | UoM | Section | Material Group
$ws.Range("A1").Value = "Material"
$ws.Range("B1").Value = "Material Description"
$ws.Range("C1").Value = "UoM"
$ws.Range("D1").Value = "Section"
$ws.Range("E1").Value = "Material Group"

# --- Column E got noticeably wider to fit the new "Material Group" header -
$ws.Range("E1").ColumnWidth = 13.6

# --- View state: zoom down a bit and move the active selection ------------
$excel.ActiveWindow.Zoom = 130
$ws.Range("B9").Select()
